$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: "Marksman" archetype ---
# Add a new "Juggler" sub-archetype under the Ranged column.
$ws.Range("B6").Value = "Juggler"
# Remove the old "Doom Sniper" tag (moved down to row 10 as "Distant Doom").
$ws.Range("L6").Value = ""

# --- Row 7: "Pugilist" archetype ---
# "God Hand" moves from column K to column S.
$ws.Range("K7").Value = ""
$ws.Range("S7").Value = "God Hand"

# --- Row 8: new "Knight" archetype ---
$ws.Range("A8").Value = "Knight"
$ws.Range("B8").Value = "ALL"
$ws.Range("F8").Value = "ALL"
$ws.Range("I8").Value = "Royal Knight"
$ws.Range("K8").Value = "White Knight"
$ws.Range("L8").Value = "Death Knight"
$ws.Range("M8").Value = "ALL"

# --- Row 9: new "Paladin" archetype ---
$ws.Range("A9").Value = "Paladin"
$ws.Range("B9").Value = "ALL"
$ws.Range("F9").Value = "ALL"
$ws.Range("G9").Value = "Holy Guard"
$ws.Range("K9").Value = "Crusader"
$ws.Range("S9").Value = "ALL"
$ws.Range("U9").Value = "Templar"

# --- Row 10: stray "Distant Doom" tag ---
$ws.Range("L10").Value = "Distant Doom"

# Update the active cell selection to match the edited workbook.
$ws.Range("B7").Select()
